$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.538.09"
$ws.Range("E2").Value = "  +2.67%  "
$ws.Range("D3").Value = "3.216.46"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.48"
$ws.Range("E5").Value = "  +2.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.56"
$ws.Range("E6").Value = "  +4.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "3.216.50"
$ws.Range("E8").Value = "  +2.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("E9").Value = "  +3.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +4.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.16"
$ws.Range("E11").Value = "  +8.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.474"
$ws.Range("E12").Value = "  +3.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("E13").Value = "  +3.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.37"
$ws.Range("E14").Value = "  +6.45%  "
$ws.Range("D15").Value = "4.755.61"
$ws.Range("E15").Value = "  +30.16%  "
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.286.64"
$ws.Range("E17").Value = "  +5.21%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.45"
$ws.Range("E18").Value = "  +4.94%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "65.188.17"
$ws.Range("E19").Value = "  +2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "485.53"
$ws.Range("E20").Value = "  +4.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.02"
$ws.Range("E21").Value = "  +5.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.770"
$ws.Range("E22").Value = "  +5.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.92"
$ws.Range("E23").Value = "  +6.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.51"
$ws.Range("E24").Value = "  +13.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.84"
$ws.Range("E25").Value = "  +6.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "83.69"
$ws.Range("E26").Value = "  +2.95%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").Value = "  +10.80%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.80"
$ws.Range("E29").Value = "  +4.42%  "
$ws.Range("E30").Value = "  +3.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.52"
$ws.Range("E31").Value = "  +7.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("E33").Value = "  +9.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.90"
$ws.Range("E34").Value = "  +7.52%  "
$ws.Range("D35").Value = "0.0₃0882"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.09"
$ws.Range("E36").Value = "  +4.56%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.40"
$ws.Range("E37").Value = "  +4.07%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.38"
$ws.Range("E38").Value = "  +6.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.48"
$ws.Range("E39").Value = "  +2.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "486.18"
$ws.Range("E40").Value = "  +10.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "52.41"
$ws.Range("E41").Value = "  +4.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.48"
$ws.Range("E42").Value = "  +8.85%  "
$ws.Range("E43").Value = "  +10.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0385"
$ws.Range("E44").Value = "  +3.94%  "
$ws.Range("D45").Value = "2.955.30"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("E46").Value = "  +3.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "39.70"
$ws.Range("E47").Value = "  +9.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.24"
$ws.Range("E48").Value = "  +5.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.35"
$ws.Range("E49").Value = "  +8.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.52"
$ws.Range("E50").Value = "  +4.83%  "
